$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so numeric-looking
# values like "592.05" are not silently converted to numbers; this
# mirrors the original file where every Price cell is an inline string.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.750.08"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.647.62"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "592.05"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "189.30"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.696"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "57.19"
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  -6.50%  "
$ws.Range("D12").Value = "0.0000270"
$ws.Range("E12").Value = "  -7.10%  "
$ws.Range("D13").Value = "10.13"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "4.232.99"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.648.08"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "18.73"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "1.10"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "67.457.26"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("D21").Value = "396.87"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").Value = "86.72"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("D27").Value = "6.04"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "3.63"
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "31.66"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "12.21"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "44.71"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "66.20"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "606.40"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "0.391"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "0.0₃0767"
$ws.Range("E40").Value = "  -14.36%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "0.0421"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E44").Value = "  -9.63%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.134"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.778.30"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").Value = "3.14"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "8.77"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").Value = "142.75"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.45"
$ws.Range("E51").Value = "  -18.03%  "

# Restore the default (unstyled) look for the Price column now that the
# values are locked in as text, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"

